$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3161.5386
$ws.Range("I64").Value = 2883.3333
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 2883.3333
$ws.Range("L64").Value = 3400
$ws.Range("M64").Value = -2635.3333
$ws.Range("N64").Value = -3896

$ws.Range("H67").Value = 3161.5386
$ws.Range("I67").Value = 2883.3333
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 2883.3333
$ws.Range("L67").Value = 3400
$ws.Range("M67").Value = -2025.3333
$ws.Range("N67").Value = -5116

$ws.Range("H74").Value = 3846.6667
$ws.Range("I74").Value = 3941.1765
$ws.Range("J74").Value = 3723.077
$ws.Range("K74").Value = 3941.1765
$ws.Range("L74").Value = 3723.077
$ws.Range("M74").Value = -3005.1765
$ws.Range("N74").Value = -5595.077

$ws.Range("H76").Value = 3088.0444
$ws.Range("I76").Value = 2999.36
$ws.Range("J76").Value = 3198.9
$ws.Range("K76").Value = 2999.36
$ws.Range("L76").Value = 3198.9
$ws.Range("M76").Value = -2684.36
$ws.Range("N76").Value = -3828.9

$ws.Range("H77").Value = 3846.6667
$ws.Range("I77").Value = 3941.1765
$ws.Range("J77").Value = 3723.077
$ws.Range("K77").Value = 19705.8825
$ws.Range("L77").Value = 18615.385
$ws.Range("M77").Value = -15025.8825
$ws.Range("N77").Value = -27975.385

$ws.Range("H79").Value = 3088.0444
$ws.Range("I79").Value = 2999.36
$ws.Range("J79").Value = 3198.9
$ws.Range("K79").Value = 2999.36
$ws.Range("L79").Value = 3198.9
$ws.Range("M79").Value = -1907.36
$ws.Range("N79").Value = -5382.9

$ws.Range("H94").Value = 3444.4443
$ws.Range("I94").Value = 3444.4443
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3444.4443
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2993.4443

$ws.Range("H120").Value = 48855
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 48855
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 48855
$ws.Range("N120").Value = -58531

$ws.Range("H128").Value = 46657.332
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 46657.332
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 46657.332
$ws.Range("N128").Value = -56617.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 39621.25
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 39621.25
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 39621.25
$ws.Range("N80").Value = -41617.25

$ws.Range("H83").Value = 39621.25
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 39621.25
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 118863.75
$ws.Range("N83").Value = -128847.75

$ws.Range("H117").Value = 50666.332
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 50666.332
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 50666.332
$ws.Range("N117").Value = -59844.332

$ws.Range("H118").Value = 50000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 50000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314

$ws.Range("H130").Value = 48429
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48429
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48429
$ws.Range("N130").Value = -58469

$ws.Range("H132").Value = 17243552
$ws.Range("I132").Value = 23810764
$ws.Range("J132").Value = 4620
$ws.Range("K132").Value = 71432292
$ws.Range("L132").Value = 13860
$ws.Range("M132").Value = -71429762
$ws.Range("N132").Value = -18920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 47473
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 47473
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 47473
$ws.Range("N117").Value = -56651

$ws.Range("H119").Value = 44796
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 44796
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 44796
$ws.Range("N119").Value = -54472

$ws.Range("H120").Value = 47761
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 47761
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 47761
$ws.Range("N120").Value = -57437

$ws.Range("H125").Value = 50780
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50780
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620

$ws.Range("H126").Value = 50772
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 50772
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652

$ws.Range("H130").Value = 49178.332
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 49178.332
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 49178.332
$ws.Range("N130").Value = -59218.332

$ws.Range("H134").Value = 2941.4067
$ws.Range("I134").Value = 1180.8
$ws.Range("J134").Value = 3844.282
$ws.Range("K134").Value = 3542.4
$ws.Range("L134").Value = 11532.846
$ws.Range("M134").Value = -1007.4
$ws.Range("N134").Value = -16602.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49883.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 49883.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 49883.5
$ws.Range("N20").Value = -50355.5

$ws.Range("H30").Value = 49883.5
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 49883.5
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 49883.5
$ws.Range("N30").Value = -50065.5

$ws.Range("H100").Value = 43436
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 43436
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 43436
$ws.Range("N100").Value = -45600

$ws.Range("H116").Value = 49822.332
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 49822.332
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 49822.332
$ws.Range("N116").Value = -59000.332

$ws.Range("H128").Value = 49883.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49883.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49883.5
$ws.Range("N128").Value = -59843.5

$ws.Range("H132").Value = 56153.42
$ws.Range("I132").Value = 1821.8889
$ws.Range("J132").Value = 178399.38
$ws.Range("K132").Value = 5465.6667
$ws.Range("L132").Value = 535198.14
$ws.Range("M132").Value = -2935.6667
$ws.Range("N132").Value = -540258.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1942.3611
$ws.Range("I132").Value = 1266.6471
$ws.Range("J132").Value = 2546.9473
$ws.Range("K132").Value = 11399.8239
$ws.Range("L132").Value = 22922.5257
$ws.Range("M132").Value = -8869.823899999999
$ws.Range("N132").Value = -27982.5257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4242.857
$ws.Range("I80").Value = 4509.091
$ws.Range("J80").Value = 3950
$ws.Range("K80").Value = 4509.091
$ws.Range("L80").Value = 3950
$ws.Range("M80").Value = -3511.091
$ws.Range("N80").Value = -5946

$ws.Range("H83").Value = 4242.857
$ws.Range("I83").Value = 4509.091
$ws.Range("J83").Value = 3950
$ws.Range("K83").Value = 22545.455
$ws.Range("L83").Value = 19750
$ws.Range("M83").Value = -17553.455
$ws.Range("N83").Value = -29734

$ws.Range("H110").Value = 49999
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 49999
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 49999
$ws.Range("N110").Value = -58179

$ws.Range("H130").Value = 53984
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 53984
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 53984
$ws.Range("N130").Value = -64024

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 46171.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 46171.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 46171.5
$ws.Range("N111").Value = -54351.5

$ws.Range("H127").Value = 50707
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50707
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50707
$ws.Range("N127").Value = -60627

$ws.Range("H130").Value = 48800
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48800
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48800
$ws.Range("N130").Value = -58840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 42676.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 42676.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 42676.75
$ws.Range("N128").Value = -52636.75

$ws.Range("H132").Value = 1439.0769
$ws.Range("I132").Value = 1263.4412
$ws.Range("J132").Value = 1770.8334
$ws.Range("K132").Value = 3790.3236
$ws.Range("L132").Value = 5312.5002
$ws.Range("M132").Value = -1260.3236
$ws.Range("N132").Value = -10372.5002

$ws.Range("H136").Value = 222996.8
$ws.Range("I136").Value = 263739.75
$ws.Range("J136").Value = 1820.8572
$ws.Range("K136").Value = 791219.25
$ws.Range("L136").Value = 5462.571599999999
$ws.Range("M136").Value = -788669.25
$ws.Range("N136").Value = -10562.5716
